$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "310.44"
Set-TextValue "E2" "-1.46%"
Set-TextValue "D3" "37.53"
Set-TextValue "E3" "-4.47%"
Set-TextValue "D4" "5.091"
Set-TextValue "E4" "-1.23%"
Set-TextValue "D5" "0.07749"
Set-TextValue "E6" "-1.46%"
Set-TextValue "E7" "-1.52%"
Set-TextValue "D8" "1.884"
Set-TextValue "E8" "-4.45%"
Set-TextValue "D9" "2.942"
Set-TextValue "E9" "-7.79%"
Set-TextValue "D10" "0.9204"
Set-TextValue "E10" "-1.80%"
Set-TextValue "D11" "0.1191"
Set-TextValue "E11" "-9.62%"
Set-TextValue "E12" "-3.53%"
Set-TextValue "D13" "0.08874"
Set-TextValue "E13" "-1.70%"
Set-TextValue "D14" "0.03390"
Set-TextValue "E14" "-2.84%"
Set-TextValue "E15" "-0.08%"
Set-TextValue "D16" "0.001377"
Set-TextValue "E16" "-2.43%"
Set-TextValue "D17" "0.005723"
Set-TextValue "E17" "-5.73%"
Set-TextValue "D18" "3.551"
Set-TextValue "E18" "-1.64%"
Set-TextValue "E19" "-1.79%"
Set-TextValue "E20" "0.44%"
Set-TextValue "D21" "0.1270"
Set-TextValue "E21" "-3.01%"
Set-TextValue "E23" "5,588.13%"
Set-TextValue "D24" "0.04396"
Set-TextValue "E24" "0.68%"
Set-TextValue "D25" "0.001212"
Set-TextValue "E25" "-2.56%"
Set-TextValue "D26" "0.004244"
Set-TextValue "E26" "-10.61%"
Set-TextValue "D27" "0.0001351"
Set-TextValue "E27" "-65.32%"
Set-TextValue "D39" "0.02116"
Set-TextValue "E39" "-5.48%"
Set-TextValue "D40" "0.04949"
Set-TextValue "E40" "-5.60%"
Set-TextValue "D41" "0.007659"
Set-TextValue "E41" "1.44%"
Set-TextValue "D42" "0.009912"
Set-TextValue "E43" "-3.88%"
Set-TextValue "E44" "-1.98%"
Set-TextValue "D45" "0.009616"
Set-TextValue "E45" "5.22%"
Set-TextValue "D46" "0.00006560"
Set-TextValue "E46" "-3.95%"
Set-TextValue "E47" "-0.16%"
Set-TextValue "E48" "0.94%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "-0.16%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "-0.16%"
